$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModelRuns")

# Insert a new row above the "2025" row (old row 73) to log the new 2023 v54 run.
# Excel copies the formatting of the row above (row 72) onto the newly
# inserted row, which matches the target styling for this entry.
$ws.Rows.Item(73).Insert()

$ws.Cells.Item(73, 1).Value = 2023
$ws.Cells.Item(73, 2).Value = "2023_TM160_IPA_54"
$ws.Cells.Item(73, 3).Value = "RTP2025_IP"
$ws.Cells.Item(73, 4).Value = "Base year"
$ws.Cells.Item(73, 6).Value = "AOC=16.21, with wrk_trn_hes=83.3"
$ws.Cells.Item(73, 7).Value = "petrale"
$ws.Cells.Item(73, 8).Value = "n/a"
$ws.Cells.Item(73, 9).Value = "current"
$ws.Cells.Item(73, 10).Value = "BlueprintNetworks_v13\net_2023_Blueprint"
$ws.Cells.Item(73, 11).Value = "model2-b"
$ws.Cells.Item(73, 12).Value = "https://app.asana.com/0/1204085012544660/1206710598691438/f"
$ws.Cells.Item(73, 13).Value = 16.21
$ws.Cells.Item(73, 14).Value = "na"
$ws.Cells.Item(73, 15).Value = "na"
$ws.Cells.Item(73, 16).Value = 1.04
$ws.Cells.Item(73, 17).Value = 0.94
$ws.Cells.Item(73, 18).Value = 83.3
$ws.Cells.Item(73, 19).Value = 0
$ws.Cells.Item(73, 20).Value = 75
$ws.Cells.Item(73, 21).Value = "2023 v54"
